# Weekly update: insert a new daily price record at the top of the data
# (row 14) for "Vega Modelo de Temuco - Camote", shifting all subsequent
# rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 14..39 down to 15..40.
$ws.Range("A14").EntireRow.Insert()

# Populate the new row 14 with this week's record.
$ws.Range("A14").Value = 10
$ws.Range("B14").Value = "Vega Modelo de Temuco"
$ws.Range("C14").Value = "La Araucanía"
$ws.Range("D14").Value = "01/18/2022"
$ws.Range("E14").Value = 9
$ws.Range("F14").Value = 100114002
$ws.Range("G14").Value = "Camote"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 50
$ws.Range("K14").Value = 20000
$ws.Range("L14").Value = 20000
$ws.Range("M14").Value = 20000
$ws.Range("N14").Value = "$/malla 20 kilos"
$ws.Range("O14").Value = "Perú"
$ws.Range("P14").Value = 1000
$ws.Range("Q14").Value = 20
$ws.Range("R14").Value = "Hortaliza"
